# Apply locality / fishing-gear / collection-date corrections cross referenced
# against the NMNH collection library.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected "Collection_Method" (and a couple of "Collection_Site" / "Collection_Date")
#     values throughout the data table ---

$ws.Range("U3").Value = "130_ft_seine"
# U3 previously used an Arial-based style (s=4); the rest of the column uses the
# Calibri-based style (s=3) that T3 already carries, so copy that formatting over.
$ws.Range("T3").Copy()
$ws.Range("U3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("T5").Value = "8_May_08"
# T5 is a brand new cell; give it the same formatting as its row neighbor S5.
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("U6").Value = "150_ft_shore_seine"
$ws.Range("U7").Value = "Dynamite"
$ws.Range("U8").Value = "150_ft_shore_seine"
$ws.Range("U11").Value = "130_ft_seine"
$ws.Range("U13").Value = "150_ft_seine"

$ws.Range("S15").Value = "Naujan_Mindoro"
$ws.Range("U15").Value = "Dynamite_25_ft_seine"
$ws.Range("U16").Value = "Dynamite_25_ft_seine"
$ws.Range("U17").Value = "Dynamite_16_ft_45_ft_seine"
$ws.Range("U19").Value = "150_ft_seine"
$ws.Range("U20").Value = "150_ft_seine"
$ws.Range("U22").Value = "250_ftm_seine"
$ws.Range("U25").Value = "Dynamite"
$ws.Range("U29").Value = "130_ft_seine"
$ws.Range("U30").Value = "Seine"
$ws.Range("U31").Value = "Dipnet_Electric_Light"
$ws.Range("U32").Value = "150_ft_seine"
$ws.Range("U34").Value = "Dipnet_Electric_Light"
$ws.Range("U35").Value = "Dynamite"

$ws.Range("S38").Value = "Port_Uson_Busuanga_Island"
$ws.Range("U38").Value = "Dynamite"

$ws.Range("U42").Value = "Dynamite"

$ws.Range("U45").Value = "130_ft_seine"
# U45 previously carried a highlighted-fill style (s=17); the corrected value
# goes back to the sheet's default (unstyled) formatting.
$ws.Range("U45").Style = "Normal"

$ws.Range("U47").Value = "130_ft_seine"
$ws.Range("U48").Value = "150_ft_seine"

# --- Previously missing Collection_Site / Collection_Date / Collection_Method
#     for rows 50-52 ---

$ws.Range("S50").Value = "Luzon_Port_San_Vicente"
$ws.Range("T50").Value = "18_Nov_08"
$ws.Range("U50").Value = "130_ft_seine"

$ws.Range("S51").Value = "Busin_Harbor_Burias_Island"
$ws.Range("T51").Value = "23_Apr_08"
$ws.Range("U51").Value = "150_ft_seine"

$ws.Range("S52").Value = "Linacapan_Island_Burias"
$ws.Range("T52").Value = "19_Dec_08"
$ws.Range("U52").Value = "130_ft_seine"

# --- View state: last-saved selection / active cell ---
$ws.Range("S50").Select()
